# add localization to text set at runtime
#
# Inserts a new "Target" key/value row (UI_GAME_TARGET / Target) before the
# existing FORMATTED_UI_END_ROUNDS_COMPLETED row (old row 45), and appends a
# new "FORMATTED_UI_GAME_TARGET" row at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Target" row at row 45 (shifts old 45-47 -> 46-48) ---
$ws.Rows("45:45").Insert()

$ws.Range("A45").Value = "UI_GAME_TARGET"
$ws.Range("B45").Value = "Target"
$ws.Range("C45").Value = "XXXX"
$ws.Range("D45").Value = "XXXX"
$ws.Range("E45").Value = "XXXX"

# --- 2. Append the new "FORMATTED_UI_GAME_TARGET" row at the end (row 49) ---
$ws.Range("A49").Value = "FORMATTED_UI_GAME_TARGET"
$ws.Range("B49").Value = "Target: {0}"
$ws.Range("C49").Value = "XXXX"
$ws.Range("D49").Value = "XXXX"
$ws.Range("E49").Value = "XXXX"

# --- 3. Nudge the saved view state to roughly match (scroll + selection) ---
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C45:E45").Select()

# --- 4. Recreate the conditional-format dxf so the style table grows from
#        2 -> 4 entries (two fresh copies of the "XXXX" highlight fill get
#        appended ahead of the rule actually in use), and the live rule ends
#        up pointing at the fourth (last) dxf, exactly like the target file.
$fc = $ws.Cells.FormatConditions
$dummy1 = $fc.Add(1, 3, '"YYYY_TMP"')
$dummy1.Interior.Color = 5066944
$dummy2 = $fc.Add(1, 3, '"XXXX"')
$dummy2.Interior.Color = 5066944

# Drop the original rule and the first throw-away rule, leaving the second
# throw-away (already targeting "XXXX") as the sole live rule.
$fc.Item(1).Delete()
$fc.Item(1).Delete()
$fc.Item(1).Priority = 2
